# Auto-generated edit script: updates D (Price) and E (Volume(1h)) columns
# for the cryptos list, per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.740.86"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "3.835.50"
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "3.835.23"
$ws.Range("E7").Value = "  -2.09%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.24"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.42%  "
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.54%  "
$ws.Range("D15").Value = "4.480.42"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").Value = "3.823.97"
$ws.Range("E16").Value = "  -2.61%  "
$ws.Range("D17").Value = "67.730.74"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.18%  "
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "466.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.46%  "
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("E24").Value = "  -4.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.82%  "
$ws.Range("E26").Value = "  -3.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("D31").Value = "3.984.79"
$ws.Range("E31").Value = "  -2.11%  "
$ws.Range("E32").Value = "  -2.39%  "
$ws.Range("E33").Value = "  -4.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.31%  "
$ws.Range("D35").Value = "3.807.99"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("E36").Value = "  -3.25%  "
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("E38").Value = "  -3.12%  "
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  -3.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "421.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("E46").Value = "  -2.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "143.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000262"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "38.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.24%  "
